$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H118").Value = 1722.7273
$ws.Range("I118").Value = 1835
$ws.Range("J118").Value = 1673.9131
$ws.Range("K118").Value = 5505
$ws.Range("L118").Value = 5021.7393
$ws.Range("M118").Value = -3848
$ws.Range("N118").Value = -8335.739300000001
$ws.Range("H132").Value = 7570.7646
$ws.Range("I132").Value = 6550.357
$ws.Range("K132").Value = 19651.071
$ws.Range("M132").Value = -17121.071
$ws.Range("H137").Value = 1549.25
$ws.Range("I137").Value = 1248.75
$ws.Range("J137").Value = 2450.75
$ws.Range("K137").Value = 3746.25
$ws.Range("L137").Value = 7352.25
$ws.Range("M137").Value = -1196.25
$ws.Range("N137").Value = -12452.25
$ws.Range("H138").Value = 2263.7246
$ws.Range("J138").Value = 2289.26
$ws.Range("L138").Value = 6867.780000000001
$ws.Range("N138").Value = -17147.78

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 353789.5
$ws.Range("I32").Value = 422127.62
$ws.Range("J32").Value = 12098.857
$ws.Range("K32").Value = 422127.62
$ws.Range("L32").Value = 12098.857
$ws.Range("M32").Value = -421840.62
$ws.Range("N32").Value = -12672.857
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 8132415.5
$ws.Range("I61").Value = 23811028
$ws.Range("J61").Value = 2764.963
$ws.Range("K61").Value = 23811028
$ws.Range("L61").Value = 2764.963
$ws.Range("M61").Value = -23810816
$ws.Range("N61").Value = -3188.963
$ws.Range("H122").Value = 1700
$ws.Range("I122").Value = 1700
$ws.Range("K122").Value = 5100
$ws.Range("M122").Value = -2650
$ws.Range("H136").Value = 8132415.5
$ws.Range("I136").Value = 23811028
$ws.Range("J136").Value = 2764.963
$ws.Range("K136").Value = 71433084
$ws.Range("L136").Value = 8294.889000000001
$ws.Range("M136").Value = -71430534
$ws.Range("N136").Value = -13394.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1423.3636
$ws.Range("I99").Value = 1405.7
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 1405.7
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = 92.29999999999995
$ws.Range("N99").Value = -4596
$ws.Range("H107").Value = 996.25
$ws.Range("I107").Value = 650.2857
$ws.Range("K107").Value = 650.2857
$ws.Range("M107").Value = 1269.7143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 18099
$ws.Range("J41").Value = 18099
$ws.Range("L41").Value = 18099
$ws.Range("N41").Value = -18955
$ws.Range("H51").Value = 16998
$ws.Range("J51").Value = 16998
$ws.Range("L51").Value = 16998
$ws.Range("N51").Value = -18470
$ws.Range("H59").Value = 23840.715
$ws.Range("J59").Value = 23840.715
$ws.Range("L59").Value = 23840.715
$ws.Range("N59").Value = -26130.715
$ws.Range("H60").Value = 13197.4
$ws.Range("J60").Value = 13197.4
$ws.Range("L60").Value = 13197.4
$ws.Range("N60").Value = -14219.4
$ws.Range("H61").Value = 16998
$ws.Range("J61").Value = 16998
$ws.Range("L61").Value = 16998
$ws.Range("N61").Value = -17694
$ws.Range("H74").Value = 28221.889
$ws.Range("J74").Value = 28221.889
$ws.Range("L74").Value = 28221.889
$ws.Range("N74").Value = -29969.889
$ws.Range("H77").Value = 28221.889
$ws.Range("J77").Value = 28221.889
$ws.Range("L77").Value = 84665.667
$ws.Range("N77").Value = -93401.667
$ws.Range("H99").Value = 1638.8518
$ws.Range("I99").Value = 983.6667
$ws.Range("J99").Value = 1720.75
$ws.Range("K99").Value = 983.6667
$ws.Range("L99").Value = 1720.75
$ws.Range("M99").Value = 514.3333
$ws.Range("N99").Value = -4716.75
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -37258
$ws.Range("H126").Value = 1638.8518
$ws.Range("I126").Value = 983.6667
$ws.Range("J126").Value = 1720.75
$ws.Range("K126").Value = 2951.0001
$ws.Range("L126").Value = 5162.25
$ws.Range("M126").Value = -481.0001000000002
$ws.Range("N126").Value = -10102.25
$ws.Range("H134").Value = 1125.3684
$ws.Range("I134").Value = 780.6667
$ws.Range("J134").Value = 1716.2858
$ws.Range("K134").Value = 2342.0001
$ws.Range("L134").Value = 5148.857400000001
$ws.Range("M134").Value = 192.9998999999998
$ws.Range("N134").Value = -10218.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2960
$ws.Range("I63").Value = 2100
$ws.Range("J63").Value = 4250
$ws.Range("K63").Value = 6300
$ws.Range("L63").Value = 12750
$ws.Range("M63").Value = -5551
$ws.Range("N63").Value = -14248
$ws.Range("H66").Value = 2960
$ws.Range("I66").Value = 2100
$ws.Range("J66").Value = 4250
$ws.Range("K66").Value = 18900
$ws.Range("L66").Value = 38250
$ws.Range("M66").Value = -15156
$ws.Range("N66").Value = -45738
$ws.Range("H68").Value = 1608.3677
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1608.3677
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4825.1031
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6447.1031
$ws.Range("H71").Value = 1608.3677
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1608.3677
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14475.3093
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22587.3093
$ws.Range("H92").Value = 762.5
$ws.Range("J92").Value = 800
$ws.Range("L92").Value = 2400
$ws.Range("N92").Value = -4896
$ws.Range("H107").Value = 1375.6104
$ws.Range("I107").Value = 271.08334
$ws.Range("J107").Value = 2345.439
$ws.Range("K107").Value = 813.2500200000001
$ws.Range("L107").Value = 7036.316999999999
$ws.Range("M107").Value = 1106.74998
$ws.Range("N107").Value = -10876.317
$ws.Range("H132").Value = 3494.7097
$ws.Range("I132").Value = 2478.074
$ws.Range("J132").Value = 4278.971
$ws.Range("K132").Value = 22302.666
$ws.Range("L132").Value = 38510.73899999999
$ws.Range("M132").Value = -19772.666
$ws.Range("N132").Value = -43570.73899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1165
$ws.Range("I97").Value = 1331.6666
$ws.Range("J97").Value = 665
$ws.Range("K97").Value = 1331.6666
$ws.Range("L97").Value = 665
$ws.Range("M97").Value = -835.6666
$ws.Range("N97").Value = -1657
$ws.Range("H122").Value = 3791.7083
$ws.Range("I122").Value = 2937.5454
$ws.Range("K122").Value = 8812.636200000001
$ws.Range("M122").Value = -6362.636200000001
$ws.Range("H130").Value = 54593.332
$ws.Range("J130").Value = 54593.332
$ws.Range("L130").Value = 54593.332
$ws.Range("N130").Value = -64633.332
$ws.Range("H132").Value = 2212.0344
$ws.Range("I132").Value = 1664.3684
$ws.Range("J132").Value = 3252.6
$ws.Range("K132").Value = 4993.1052
$ws.Range("L132").Value = 9757.799999999999
$ws.Range("M132").Value = -2463.1052
$ws.Range("N132").Value = -14817.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45457010
$ws.Range("I7").Value = 66668930
$ws.Range("K7").Value = 66668930
$ws.Range("M7").Value = -66668818
$ws.Range("H40").Value = 47621656
$ws.Range("I40").Value = 55557684
$ws.Range("J40").Value = 5500
$ws.Range("K40").Value = 55557684
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = -55557548
$ws.Range("N40").Value = -5772
$ws.Range("H122").Value = 3358.5366
$ws.Range("I122").Value = 2357.1428
$ws.Range("J122").Value = 4410
$ws.Range("K122").Value = 7071.428400000001
$ws.Range("L122").Value = 13230
$ws.Range("M122").Value = -4621.428400000001
$ws.Range("N122").Value = -18130
$ws.Range("H126").Value = 45457010
$ws.Range("I126").Value = 66668930
$ws.Range("K126").Value = 200006790
$ws.Range("M126").Value = -200004320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4901.8184
$ws.Range("I81").Value = 5092
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 10184
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -9123
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 4901.8184
$ws.Range("I84").Value = 5092
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 50920
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -45616
$ws.Range("N84").Value = -40608
$ws.Range("H100").Value = 922
$ws.Range("I100").Value = 708
$ws.Range("J100").Value = 1136
$ws.Range("K100").Value = 1416
$ws.Range("L100").Value = 2272
$ws.Range("M100").Value = -875
$ws.Range("N100").Value = -3354
$ws.Range("H122").Value = 1968.75
$ws.Range("I122").Value = 1581.8182
$ws.Range("K122").Value = 4745.4546
$ws.Range("M122").Value = -2295.4546
$ws.Range("H126").Value = 1376.6428
$ws.Range("I126").Value = 995.5
$ws.Range("J126").Value = 1884.8334
$ws.Range("K126").Value = 2986.5
$ws.Range("L126").Value = 5654.5002
$ws.Range("M126").Value = -516.5
$ws.Range("N126").Value = -10594.5002
